$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "Earn an Acoustic Engineering Degree | Acoustical Schools"
$ws.Cells.Item(2,2).Value = "https://educatingengineers.com/degrees/acoustic-engineering"
$ws.Cells.Item(2,3).Value = 68
$ws.Cells.Item(2,4).Value = 4
$ws.Cells.Item(2,5).Value = 8
$ws.Cells.Item(2,6).Value = ""

# Row 3
$ws.Cells.Item(3,1).Value = "Compare Industrial Engineering Courses and Degree Programs"
$ws.Cells.Item(3,2).Value = "https://educatingengineers.com/degrees/industrial-engineering"
$ws.Cells.Item(3,3).Value = 55
$ws.Cells.Item(3,4).Value = 7
$ws.Cells.Item(3,5).Value = 6
$ws.Cells.Item(3,6).Value = ""

# Row 4
$ws.Cells.Item(4,1).Value = "Nuclear Engineering Schools and Degrees | EducatingEngineers.com"
$ws.Cells.Item(4,2).Value = "https://educatingengineers.com/degrees/nuclear-engineering"
$ws.Cells.Item(4,3).Value = 52
$ws.Cells.Item(4,4).Value = 5
$ws.Cells.Item(4,5).Value = 5
$ws.Cells.Item(4,6).Value = ""

# Row 5
$ws.Cells.Item(5,1).Value = "Civil Engineer Jobs and Careers | EducatingEngineers.com"
$ws.Cells.Item(5,2).Value = "https://educatingengineers.com/careers/civil-engineer"
$ws.Cells.Item(5,3).Value = 134
$ws.Cells.Item(5,4).Value = 37
$ws.Cells.Item(5,5).Value = 24
$ws.Cells.Item(5,6).Value = ""

# Row 6
$ws.Cells.Item(6,1).Value = "8 Best Engineering Jobs | Best Jobs Rankings | US News Careers"
$ws.Cells.Item(6,2).Value = "https://money.usnews.com/careers/best-jobs/rankings/best-engineering-jobs"
$ws.Cells.Item(6,3).Value = 661
$ws.Cells.Item(6,4).Value = 4
$ws.Cells.Item(6,5).Value = 55
$ws.Cells.Item(6,6).Value = ""

# Row 7
$ws.Cells.Item(7,1).Value = "List of Engineering Career Options with Job Descriptions | EducatingEngineers.com"
$ws.Cells.Item(7,2).Value = "https://educatingengineers.com/career-specialties"
$ws.Cells.Item(7,3).Value = 95
$ws.Cells.Item(7,4).Value = 18
$ws.Cells.Item(7,5).Value = 6
$ws.Cells.Item(7,6).Value = ""

# Row 8
$ws.Cells.Item(8,1).Value = "Engineering Careers: Options, Job Titles, and Descriptions"
$ws.Cells.Item(8,2).Value = "https://www.thebalancecareers.com/engineering-job-titles-2061493"
$ws.Cells.Item(8,3).Value = 339
$ws.Cells.Item(8,4).Value = 13
$ws.Cells.Item(8,5).Value = 27
$ws.Cells.Item(8,6).Value = ""

# Row 9
$ws.Cells.Item(9,1).Value = "Top-paying jobs are in engineering "
$ws.Cells.Item(9,2).Value = "https://money.cnn.com/2013/04/25/news/economy/engineering-best-paid-jobs/index.html?sa=X&ved=2ahUKEwiS28uDwZzmAhUJPa0KHRRTAPQQ9QF6BAgLEAI"
$ws.Cells.Item(9,3).Value = 26
$ws.Cells.Item(9,4).Value = 4
$ws.Cells.Item(9,5).Value = 35
$ws.Cells.Item(9,6).Value = ""

# Row 10
$ws.Cells.Item(10,1).Value = "The Best Engineering Jobs for Engineers | ENGINEERING.com"
$ws.Cells.Item(10,2).Value = "https://www.engineering.com/jobs/"
$ws.Cells.Item(10,3).Value = 362
$ws.Cells.Item(10,4).Value = 2
$ws.Cells.Item(10,5).Value = 7
$ws.Cells.Item(10,6).Value = "support@engineering.com`n"
$ws.Rows.Item(10).AutoFit()

# Row 11 (new)
$ws.Cells.Item(11,1).Value = "The 6 Highest Paid Engineering Jobs"
$ws.Cells.Item(11,2).Value = "https://typesofengineeringdegrees.org/highest-paid-engineering-jobs/"
$ws.Cells.Item(11,3).Value = 95
$ws.Cells.Item(11,4).Value = 37
$ws.Cells.Item(11,5).Value = 13
$ws.Cells.Item(11,6).Value = ""
